$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @(2024, "Richard",   1, 11, 0, 11, 32850, 30, 10, "Active", 366),
    @(2024, "Andy",      2,  9, 0,  9, 29350, 50, 30, "Active", 349),
    @(2024, "Matt",      3,  8, 0,  8, 35650, 10,-10, "Active", 362),
    @(2024, "Jon",       4,  7, 0,  7, 18350,  0,-20, "Active", 357),
    @(2024, "Prashant",  5,  6, 0,  6, 23100, 50, 40, "Active", 365),
    @(2024, "Maisy",     6,  5, 0,  5, 21800, 20,  0, "Active", 360),
    @(2024, "Pepe",      7,  5, 0,  5, 13900,  0,-20, "Active", 364),
    @(2024, "Anthony",   8,  5, 0,  5, 11750,  0,-10, "Active", 350),
    @(2024, "Mark",      9,  4, 0,  4, 12050,  0,-20, "Active", 361)
)

$startRow = 191
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]
    $ws.Cells.Item($row, 6).Value = $vals[5]
    $ws.Cells.Item($row, 7).Value = $vals[6]
    $ws.Cells.Item($row, 8).Value = $vals[7]
    $ws.Cells.Item($row, 9).Value = $vals[8]
    $ws.Cells.Item($row, 10).Value = $vals[9]
    $ws.Cells.Item($row, 11).Value = $vals[10]
}

[void]$ws.Range("A1:K199").Select()
